# Add "Source" / "State and Territory governments" row to the Description
# sheet, then make Description the active/selected sheet with A5:B5
# selected (mirrors the author switching to the Description tab after
# adding the new sourcing row for all NPAs).

$wb = $excel.ActiveWorkbook

$wsDescription = $wb.Worksheets.Item("Description")

# New row 5: Source -> State and Territory governments
$wsDescription.Range("A5").Value = "Source"
$wsDescription.Range("B5").Value = "State and Territory governments"

# Make Description the active sheet and select the newly added row.
$wsDescription.Activate() | Out-Null
$wsDescription.Range("A5:B5").Select() | Out-Null
